$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 2806.3
$ws.Range("B3").Value = 2695.5
$ws.Range("C3").Value = 2768.2
$ws.Range("C4").Value = 2688.8
$ws.Range("C5").Value = 2622.6
$ws.Range("C6").Value = 2578.9
$ws.Range("C9").Value = 2189.2
$ws.Range("C12").Value = 2315
$ws.Range("C15").Value = 2679.2
$ws.Range("C18").Value = 3432.8
